$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.909.63"
$ws.Range("E2").Value = "  +1.28%  "
$ws.Range("D3").Value = "1.644.00"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.55"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.52"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.09%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0616"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.70%  "
$ws.Range("E11").Value = "  -1.49%  "
$ws.Range("D12").Value = "1.877.00"
$ws.Range("E12").Value = "  +1.28%  "
$ws.Range("D13").Value = "1.649.56"
$ws.Range("E13").Value = "  +1.70%  "
$ws.Range("E14").Value = "  +4.22%  "
$ws.Range("E15").Value = "  +0.42%  "
$ws.Range("E16").Value = "  +0.91%  "
$ws.Range("D17").Value = "27.896.29"
$ws.Range("E17").Value = "  +1.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.66"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.49%  "
$ws.Range("D19").Value = "0.0₃0724"
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("E20").Value = "  +0.84%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.83"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.86%  "
$ws.Range("E23").Value = "  +1.28%  "
$ws.Range("E24").Value = "  +2.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.56"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.92"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.51%  "
$ws.Range("E27").Value = "  +0.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.72"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  +1.09%  "
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("D33").Value = "1.427.31"
$ws.Range("E33").Value = "  -2.78%  "
$ws.Range("E34").Value = "  +0.40%  "
$ws.Range("E35").Value = "  +1.36%  "
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.885"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.38%  "
$ws.Range("E38").Value = "  +0.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.928"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.557"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.26%  "
$ws.Range("E41").Value = "  +1.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "68.63"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.24%  "
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("E45").Value = "  +2.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.80"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.52%  "
$ws.Range("E47").Value = "  +0.15%  "
$ws.Range("D48").Value = "1.785.60"
$ws.Range("E48").Value = "  +1.15%  "
$ws.Range("E49").Value = "  +1.87%  "
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0506"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.51%  "
